$d = $word.ActiveDocument

# 1) Remove the "Meta description" paragraph that currently sits right
#    after the H1 title ("Play Christmas Cash Pots Free | Festive Slot Game").
$metaOld = 'Meta description: Get into the holiday spirit and play Christmas Cash Pots for free. Review of the immersive festive slot game with a bonus feature and high RTP.'
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $metaOld) {
        $p.Range.Delete() | Out-Null
        break
    }
}

# 2) Insert a new bold paragraph ("Play Christmas Cash Pots Free | Festive
#    Slot Game") right before the final "Prompt for DALLE" paragraph.
$count = $d.Paragraphs.Count
$secondLast = $d.Paragraphs.Item($count - 1)
$secondLast.Range.InsertParagraphAfter() | Out-Null

$count = $d.Paragraphs.Count
$newPara = $d.Paragraphs.Item($count - 1)
$xmlFrag = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Christmas Cash Pots Free | Festive Slot Game</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$newPara.Range.InsertXML($xmlFrag) | Out-Null

# 3) Replace the old DALLE image-prompt text (now the last paragraph) with
#    the review's meta description, keeping its italic run formatting.
$old = 'Prompt for DALLE: Create a feature image for the game "Christmas Cash Pots" in a cartoon style. The image should showcase a happy Maya warrior wearing glasses.'
$new = 'Get into the holiday spirit and play Christmas Cash Pots for free. Review of the immersive festive slot game with a bonus feature and high RTP.'
$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
